$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.070.04"
$ws.Range("E2").Value = "  -2.90%  "

$ws.Range("D3").Value = "1.718.12"
$ws.Range("E3").Value = "  -2.97%  "

$ws.Range("D4").Value = "'1.006"
$ws.Range("E4").Value = "  +0.21%  "

$ws.Range("D5").Value = "'316.92"
$ws.Range("E5").Value = "  -3.28%  "

$ws.Range("D6").Value = "'1.005"
$ws.Range("E6").Value = "  +0.26%  "

$ws.Range("D7").Value = "'0.4644"
$ws.Range("E7").Value = "  +3.52%  "

$ws.Range("D8").Value = "'0.3435"
$ws.Range("E8").Value = "  -3.80%  "

$ws.Range("D9").Value = "'42.43"
$ws.Range("E9").Value = "  +0.61%  "

$ws.Range("D10").Value = "'0.07287"
$ws.Range("E10").Value = "  -2.38%  "

$ws.Range("D11").Value = "'1.049"
$ws.Range("E11").Value = "  -4.24%  "

$ws.Range("D12").Value = "'1.004"
$ws.Range("E12").Value = "  +0.15%  "

$ws.Range("D13").Value = "'19.88"
$ws.Range("E13").Value = "  -5.01%  "

$ws.Range("D14").Value = "'5.868"
$ws.Range("E14").Value = "  -3.05%  "

$ws.Range("D15").Value = "1.727.73"
$ws.Range("E15").Value = "  -2.49%  "

$ws.Range("E16").Value = "  -4.44%  "

$ws.Range("D17").Value = "'89.64"
$ws.Range("E17").Value = "  -3.46%  "

$ws.Range("D18").Value = "'0.00001046"
$ws.Range("E18").Value = "  -1.41%  "

$ws.Range("D19").Value = "'0.06288"
$ws.Range("E19").Value = "  -2.07%  "

$ws.Range("D20").Value = "'1.006"
$ws.Range("E20").Value = "  +0.40%  "

$ws.Range("E21").Value = "  -4.31%  "

$ws.Range("E22").Value = "  -3.55%  "

$ws.Range("D23").Value = "27.151.42"
$ws.Range("E23").Value = "  -2.70%  "

$ws.Range("D24").Value = "'10.80"
$ws.Range("E24").Value = "  -4.73%  "

$ws.Range("D25").Value = "'2.150"
$ws.Range("E25").Value = "  +1.67%  "

$ws.Range("D26").Value = "'156.61"
$ws.Range("E26").Value = "  -3.73%  "

$ws.Range("D27").Value = "'19.46"
$ws.Range("E27").Value = "  -4.00%  "

$ws.Range("D28").Value = "1.931.55"
$ws.Range("E28").Value = "  -2.25%  "

$ws.Range("D29").Value = "'2.137"
$ws.Range("E29").Value = "  -2.78%  "

$ws.Range("D30").Value = "'119.12"
$ws.Range("E30").Value = "  -5.19%  "

$ws.Range("D31").Value = "'1.020"
$ws.Range("E31").Value = "  -7.48%  "

$ws.Range("D32").Value = "'0.09091"
$ws.Range("E32").Value = "  -0.81%  "

$ws.Range("D33").Value = "'3.603"
$ws.Range("E33").Value = "  -0.95%  "

$ws.Range("D34").Value = "'5.321"
$ws.Range("E34").Value = "  -4.70%  "

$ws.Range("D35").Value = "'0.02200"
$ws.Range("E35").Value = "  -4.20%  "

$ws.Range("D36").Value = "'11.13"
$ws.Range("E36").Value = "  -6.30%  "

$ws.Range("D37").Value = "'0.05832"

$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").Value = "'4.761"
$ws.Range("E38").Value = "  -4.22%  "

$ws.Range("B39").Value = "Algorand"
$ws.Range("C39").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D39").Value = "'0.1989"
$ws.Range("E39").Value = "  -5.18%  "

$ws.Range("D40").Value = "'1.407"
$ws.Range("E40").Value = "  +1.03%  "

$ws.Range("D41").Value = "'0.5961"
$ws.Range("E41").Value = "  -5.99%  "

$ws.Range("D42").Value = "'1.132"
$ws.Range("E42").Value = "  -4.48%  "

$ws.Range("D43").Value = "'7.479"
$ws.Range("E43").Value = "  -5.66%  "

$ws.Range("D44").Value = "'3.639"
$ws.Range("E44").Value = "  -2.72%  "

$ws.Range("D45").Value = "'12.63"
$ws.Range("E45").Value = "  -4.68%  "

$ws.Range("D46").Value = "'0.5605"
$ws.Range("E46").Value = "  -4.63%  "

$ws.Range("D47").Value = "'119.32"
$ws.Range("E47").Value = "  -2.76%  "

$ws.Range("D48").Value = "'1.862"
$ws.Range("E48").Value = "  -5.03%  "

$ws.Range("D49").Value = "'0.06668"
$ws.Range("E49").Value = "  -3.63%  "

$ws.Range("D50").Value = "'1.087"
$ws.Range("E50").Value = "  -4.72%  "

$ws.Range("D51").Value = "'1.004"
$ws.Range("E51").Value = "  +0.30%  "
